$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF").
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting used by the existing header cells (bold font, thin
# border, centered/top aligned) by copying H1's format onto the new headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I ("I0") and J ("IF") for rows 2-19.
$data = @(
    @(5, 5),
    @(7, 7),
    @(6, 7),
    @(3, 4),
    @(5, 5),
    @(6, 7),
    @(6, 7),
    @(7, 7),
    @(7, 9),
    @(6, 7),
    @(7, 8),
    @(8, 9),
    @(7, 8),
    @(1, 5),
    @(1, 4),
    @(1, 6),
    @(1, 3),
    @(7, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $pair = $data[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
